$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the cells whose value actually changes ---
$ws.Range("A2").Value = "UJLnj648"
$ws.Range("B2").Value = 23110820
$ws.Range("C2").Value = "usomuin80"
$ws.Range("D2").Value = "u#A74V!n"
$ws.Range("F2").Value = "XCqBUXBb"
$ws.Range("G2").Value = "NJae"

# --- Row 3: only the cells whose value actually changes ---
$ws.Range("A3").Value = "opEDi265"
$ws.Range("B3").Value = 23110819
$ws.Range("C3").Value = "pptaghb76"
$ws.Range("D3").Value = "eTC7u&2$"
$ws.Range("F3").Value = "ENpBoHRd"
$ws.Range("G3").Value = "dFMq"

# --- Row 4: brand-new row, same layout/style as row 3 ---
$ws.Range("A4").Value = "GaiRe814"
$ws.Range("B4").Value = 23110818
$ws.Range("C4").Value = "fafepfu81"
$ws.Range("D4").Value = "c7C9N!#p"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "OBLxNjot"
$ws.Range("G4").Value = "gZgz"
$ws.Range("H4").Value = "Candidate"

# Match styling of the row above for the new row
$ws.Range("A4:H4").Style = $ws.Range("A3:H3").Style

# Keep the selection/used-range in sync with the new data
$ws.Range("A1:H4").Select()
